$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.264.16"
$ws.Range("E2").Value = "  +0.34%  "

$ws.Range("D3").Value = "1.689.18"
$ws.Range("E3").Value = "  +1.32%  "

$ws.Range("E4").Value = "  +0.16%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "218.75"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.69%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5249"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +2.99%  "

$ws.Range("E7").Value = "  +0.12%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2695"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +1.81%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "22.02"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +2.26%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07467"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +1.44%  "

$ws.Range("D12").Value = "1.705.81"
$ws.Range("E12").Value = "  +2.21%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.556"
$ws.Range("D13").ClearFormats()

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5851"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +1.39%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.000008505"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -0.10%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.56"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -0.50%  "

$ws.Range("D17").Value = "26.314.80"
$ws.Range("E17").Value = "  +0.31%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.963"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.70%  "

$ws.Range("E19").Value = "  +0.12%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.86"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.35%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "189.12"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.25%  "

$ws.Range("E22").Value = "  +0.66%  "

$ws.Range("E23").Value = "  +0.08%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "144.77"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +1.23%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "7.667"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.39%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1234"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +5.42%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.84"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +0.97%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.06691"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +15.43%  "

$ws.Range("E29").Value = "  +5.67%  "

$ws.Range("E31").Value = "  +2.48%  "

$ws.Range("E32").Value = "  +1.19%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.669"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +1.41%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.029"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +2.37%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6218"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +3.92%  "

$ws.Range("E36").Value = "  +1.54%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.706"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +2.60%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.312"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +5.31%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01621"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +0.90%  "

$ws.Range("D40").Value = "1.102.78"
$ws.Range("E40").Value = "  +1.98%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8852"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +3.19%  "

$ws.Range("E42").Value = "  +0.76%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "101.29"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +1.46%  "

$ws.Range("D44").Value = "1.836.51"
$ws.Range("E44").Value = "  +1.08%  "

$ws.Range("E45").Value = "  +2.90%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "56.80"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +2.00%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.180"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +1.41%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.006"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.01%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05263"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +1.62%  "

$ws.Range("E50").Value = "  +0.12%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.059"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +3.29%  "
